$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: populate new rows 3-7 with the correct style "shapes" first, by
# copying formats from existing rows that already carry the right combination
# of font / fill / number-format, so the engine reuses existing cellXfs
# entries instead of minting new ones. ---

# Row 3 shares Row 2's style shape (E/F have the "no-fill / fill" blank-amount
# pattern). Copy it BEFORE H2's own style changes below.
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122)

# Row 4 and Row 6 share Row 1's style shape (E/F have the "amount / fill"
# pattern).
$ws.Range("A1:I1").Copy()
$ws.Range("A4:I4").PasteSpecial(-4122)
$ws.Range("A1:I1").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)

# Row 5 also shares Row 1's style shape for columns A-H.
$ws.Range("A1:H1").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)

# Row 7 shares Row 2/3's style shape.
$ws.Range("A2:I2").Copy()
$ws.Range("A7:I7").PasteSpecial(-4122)

# J4 gets the same "plain text" style used by A/C/G/H columns.
$ws.Range("A1").Copy()
$ws.Range("J4").PasteSpecial(-4122)

# I5 keeps the "highlighted / filled" look used by column I elsewhere.
$ws.Range("I1").Copy()
$ws.Range("I5").PasteSpecial(-4122)

# J5 is a blank note cell using the bigger "Calibri" font (same font as the
# amount columns use when highlighted) but with a plain General number
# format.
$ws.Range("E2").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Cells.Item(5, 10).NumberFormat = "General"

$ws.Cells.Item(1, 1).Copy()   # clear clipboard/marching ants

# --- Step 2: cell values ---

# Row 2: only H2 changes, from "153 Orange" to "Financial Service"
$ws.Cells.Item(2, 8).Value2 = "Financial Service"

# Row 3
$ws.Cells.Item(3, 1).Value2 = "BoA-7797"
$ws.Cells.Item(3, 2).Value2 = 44635
$ws.Cells.Item(3, 3).Value2 = "City of Philadelphia Bill Payment"
$ws.Cells.Item(3, 4).Value2 = 737
$ws.Cells.Item(3, 5).ClearContents()
$ws.Cells.Item(3, 7).Value2 = "Taxes"
$ws.Cells.Item(3, 8).Value2 = "332 Robbins"
$ws.Cells.Item(3, 9).Value2 = "Taxes"

# Row 4
$ws.Cells.Item(4, 1).Value2 = "BoA-7797"
$ws.Cells.Item(4, 2).Value2 = 44706
$ws.Cells.Item(4, 3).Value2 = "ZelleTony for Beacon clean sewer pipe etc."
$ws.Cells.Item(4, 4).Value2 = 3500
$ws.Cells.Item(4, 5).Value2 = 2200
$ws.Cells.Item(4, 7).Value2 = "Repairs"
$ws.Cells.Item(4, 8).Value2 = "207 Beacon"
$ws.Cells.Item(4, 9).Value2 = "Repairs"
$ws.Cells.Item(4, 10).Value2 = "Steven "

# Row 5
$ws.Cells.Item(5, 1).Value2 = "BoA-7797"
$ws.Cells.Item(5, 2).Value2 = 44706
$ws.Cells.Item(5, 3).Value2 = "ZelleTony for Beacon clean sewer pipe etc."
$ws.Cells.Item(5, 4).Value2 = 3500
$ws.Cells.Item(5, 5).Value2 = 480
$ws.Cells.Item(5, 7).Value2 = "Cleaning & maintenance"
$ws.Cells.Item(5, 8).Value2 = "207 Beacon"
$ws.Cells.Item(5, 9).Value2 = "Cleaning & maintenance"

# Row 6
$ws.Cells.Item(6, 1).Value2 = "BoA-8211"
$ws.Cells.Item(6, 2).Value2 = 44662
$ws.Cells.Item(6, 3).Value2 = "RAPIN FRITURA RESTAURANT NEWARK NJ"
$ws.Cells.Item(6, 4).Value2 = 58.45
$ws.Cells.Item(6, 5).Value2 = 29.23
$ws.Cells.Item(6, 7).Value2 = "Meal - C"
$ws.Cells.Item(6, 8).Value2 = "Financial Service"
$ws.Cells.Item(6, 9).Value2 = "Meal - C"

# Row 7
$ws.Cells.Item(7, 1).Value2 = "BoA-8211"
$ws.Cells.Item(7, 2).Value2 = 44858
$ws.Cells.Item(7, 3).Value2 = "OPTIMUM 7875 973-230-6046 NY"
$ws.Cells.Item(7, 4).Value2 = 50.65
$ws.Cells.Item(7, 5).ClearContents()
$ws.Cells.Item(7, 7).Value2 = "Phone/Cable/Internet"
$ws.Cells.Item(7, 8).Value2 = "207 Beacon"
$ws.Cells.Item(7, 9).Value2 = "Other"

# --- Step 3: formulas for column F, rows 3-7 (mirrors the F1/F2 formula) ---
for ($r = 3; $r -le 7; $r++) {
    $formula = '=if(And(G' + $r + '<>"",H' + $r + '<>""),if(E' + $r + '<>"",E' + $r + ',D' + $r + '),)'
    $ws.Cells.Item($r, 6).Formula = $formula
}

# --- Step 4: extend the G/H data validation dropdowns down to row 7 ---
$ws.Range("H1:H7").Validation.Delete()
$ws.Range("G1:G7").Validation.Delete()
$ws.Range("H1:H7").Validation.Add(3, 1, 1, "#REF!")
$ws.Range("G1:G7").Validation.Add(3, 1, 1, "#REF!")

$wb.Save()
